$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the extra trailing rows (rows 18-20 in the old layout) so only
# 16 data rows (+1 header) remain.
$ws.Range("A18:D20").EntireRow.Delete() | Out-Null

# Add the new "ITI" header in column D.
$ws.Range("D1").Value = "ITI"

# Update ConditionType (C) and ITI (D) values for each trial row.
$data = @(
    @(2, 7),
    @(4, 9),
    @(1, 6),
    @(4, 7),
    @(4, 6),
    @(3, 6),
    @(1, 6),
    @(3, 8),
    @(2, 6),
    @(2, 6),
    @(1, 6),
    @(1, 7),
    @(2, 8),
    @(3, 7),
    @(4, 7),
    @(3, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $data[$i][0]
    $ws.Cells.Item($row, 4).Value = $data[$i][1]
}

$ws.Range("E22").Select() | Out-Null
